$wb = $excel.ActiveWorkbook

# Sheet "Hoja1": update the daily conversion text in A1
$ws1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.54 = 9595.42 pesos`n✅ 9595.42 pesos = 2.54 = 945.26 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $newText

# Sheet "tasas": update the four rate figures
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 393
$ws2.Range("O10").Value = 3771
$ws2.Range("N12").Value = 3780
$ws2.Range("O12").Value = 372.372
